# Adds two scikit-learn style classification-report tables to the
# "Hypertension / Summary tables" sheet, in the block of small summary
# tables that live to the right of the main descriptive-stats tables
# (columns N:R), directly below the existing Sleep_Duration/BMI table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Table 1: rows 35-41  (classification report #1)
# ---------------------------------------------------------------------

# Row 35: blank spacer row (matches the blank spacer style used elsewhere,
# e.g. row 34 / row 3, with R getting the borderless-blank style used at R9)
$ws.Range("C34").Copy()
$ws.Range("N35:Q35").PasteSpecial(-4122)
$ws.Range("R9").Copy()
$ws.Range("R35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 36: header row (Class | Precision | Recall | F1-score | Support)
$ws.Range("N30:P30").Copy()
$ws.Range("N36:P36").PasteSpecial(-4122)
$ws.Range("P30").Copy()
$ws.Range("Q36:R36").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("N36").Value = "Class"
$ws.Range("O36").Value = "Precision"
$ws.Range("P36").Value = "Recall"
$ws.Range("Q36").Value = "F1-score"
$ws.Range("R36").Value = "Support"

# Row 37: "No" class metrics
$ws.Range("N11").Copy()
$ws.Range("N37").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("O37:R37").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("N37").Value = "No"
$ws.Range("O37").Value = 0.83
$ws.Range("P37").Value = 0.87
$ws.Range("Q37").Value = 0.85
$ws.Range("R37").Value = 191

# Row 38: "Yes" class metrics
$ws.Range("N11").Copy()
$ws.Range("N38").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("O38:R38").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("N38").Value = "Yes"
$ws.Range("O38").Value = 0.87
$ws.Range("P38").Value = 0.83
$ws.Range("Q38").Value = 0.85
$ws.Range("R38").Value = 206

# Row 39: "Accuracy" row (only one metric column is populated; O/P show a
# placeholder dash, Q/R use a distinct bold centred style)
$ws.Range("N11").Copy()
$ws.Range("N39").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("O39:P39").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("N39").Value = "Accuracy"
$ws.Range("O39").Value = "—"
$ws.Range("P39").Value = "—"
$ws.Range("Q39").Value = 0.85
$ws.Range("Q39").Font.Bold = $true
$ws.Range("Q39").HorizontalAlignment = -4108
$ws.Range("Q39").VerticalAlignment = -4108
$ws.Range("Q39").WrapText = $true
$ws.Range("Q39").Copy()
$ws.Range("R39").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("R39").Value = 397

# Row 40: "Macro Avg" row
$ws.Range("N11").Copy()
$ws.Range("N40").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("O40:R40").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("N40").Value = "Macro Avg"
$ws.Range("O40").Value = 0.85
$ws.Range("P40").Value = 0.85
$ws.Range("Q40").Value = 0.85
$ws.Range("R40").Value = 397

# Row 41: "Weighted Avg" row
$ws.Range("N12").Copy()
$ws.Range("N41").PasteSpecial(-4122)
$ws.Range("D9").Copy()
$ws.Range("O41:R41").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("N41").Value = "Weighted Avg"
$ws.Range("O41").Value = 0.85
$ws.Range("P41").Value = 0.85
$ws.Range("Q41").Value = 0.85
$ws.Range("R41").Value = 397

# ---------------------------------------------------------------------
# Table 2: rows 43-49  (classification report #2)
# ---------------------------------------------------------------------

# Row 43: blank spacer row
$ws.Range("C34").Copy()
$ws.Range("N43:Q43").PasteSpecial(-4122)
$ws.Range("R9").Copy()
$ws.Range("R43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 44: header row
$ws.Range("N30:P30").Copy()
$ws.Range("N44:P44").PasteSpecial(-4122)
$ws.Range("P30").Copy()
$ws.Range("Q44:R44").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("N44").Value = "Class"
$ws.Range("O44").Value = "Precision"
$ws.Range("P44").Value = "Recall"
$ws.Range("Q44").Value = "F1-score"
$ws.Range("R44").Value = "Support"

# Row 45: "No" class metrics
$ws.Range("N16").Copy()
$ws.Range("N45").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("O45:R45").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("N45").Value = "No"
$ws.Range("O45").Value = 0.84
$ws.Range("P45").Value = 0.83
$ws.Range("Q45").Value = 0.83
$ws.Range("R45").Value = 191

# Row 46: "Yes" class metrics
$ws.Range("N16").Copy()
$ws.Range("N46").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("O46:R46").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("N46").Value = "Yes"
$ws.Range("O46").Value = 0.84
$ws.Range("P46").Value = 0.85
$ws.Range("Q46").Value = 0.85
$ws.Range("R46").Value = 206

# Row 47: "Accuracy" row
$ws.Range("N11").Copy()
$ws.Range("N47").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("O47:P47").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("N47").Value = "Accuracy"
$ws.Range("O47").Value = "—"
$ws.Range("P47").Value = "—"
$ws.Range("Q39").Copy()
$ws.Range("Q47").PasteSpecial(-4122)
$ws.Range("R47").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Q47").Value = 0.84
$ws.Range("R47").Value = 397

# Row 48: "Macro Avg" row
$ws.Range("N11").Copy()
$ws.Range("N48").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("O48:R48").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("N48").Value = "Macro Avg"
$ws.Range("O48").Value = 0.84
$ws.Range("P48").Value = 0.84
$ws.Range("Q48").Value = 0.84
$ws.Range("R48").Value = 397

# Row 49: "Weighted Avg" row
$ws.Range("N12").Copy()
$ws.Range("N49").PasteSpecial(-4122)
$ws.Range("D9").Copy()
$ws.Range("O49:R49").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("N49").Value = "Weighted Avg"
$ws.Range("O49").Value = 0.84
$ws.Range("P49").Value = 0.84
$ws.Range("Q49").Value = 0.84
$ws.Range("R49").Value = 397

# ---------------------------------------------------------------------
# View state: scroll position + active selection, matching the author's
# final cursor position when they finished editing.
# ---------------------------------------------------------------------
$ws.Range("L42").Select()
